$p = $ppt.ActivePresentation

# The deck originally opened with a generic "RL Seminar" title slide
# (old slide 1). That slide is being removed entirely; slide 2 (the
# "Value Function Approximation" title slide) becomes the new slide 1,
# and every other slide shifts up by one position.
$p.Slides.Item(1).Delete()
